# Refresh "cryptos" list (Price / Volume(1h) columns) with the latest scrape.
# Values that look like plain numbers are written with a leading apostrophe so
# Excel keeps them as literal text (matching the original inlineStr cells)
# instead of silently reparsing them as numbers (which would drop things like
# trailing zeros, e.g. "8.130" -> 8.13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.619.50'
$ws.Range('E2').Value = '  +3.21%  '
$ws.Range('D3').Value = '1.697.91'
$ws.Range('E3').Value = '  +2.02%  '
$ws.Range('D4').Value = '''0.9991'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''316.12'
$ws.Range('E5').Value = '  +1.82%  '
$ws.Range('D6').Value = '''0.9996'
$ws.Range('D7').Value = '''0.3942'
$ws.Range('E7').Value = '  +1.45%  '
$ws.Range('D8').Value = '''0.4013'
$ws.Range('E8').Value = '  +0.92%  '
$ws.Range('D9').Value = '''1.536'
$ws.Range('E9').Value = '  +3.86%  '
$ws.Range('D10').Value = '''0.9999'
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('D11').Value = '''53.46'
$ws.Range('E11').Value = '  +3.98%  '
$ws.Range('D12').Value = '''0.08762'
$ws.Range('E12').Value = '  +0.94%  '
$ws.Range('D13').Value = '''7.232'
$ws.Range('E13').Value = '  +7.24%  '
$ws.Range('D14').Value = '''23.24'
$ws.Range('E14').Value = '  +2.41%  '
$ws.Range('D15').Value = '''8.130'
$ws.Range('E15').Value = '  +11.10%  '
$ws.Range('D16').Value = '''0.00001317'
$ws.Range('E16').Value = '  +0.65%  '
$ws.Range('D17').Value = '1.695.26'
$ws.Range('E17').Value = '  +2.01%  '
$ws.Range('D18').Value = '''99.75'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').Value = '''0.07075'
$ws.Range('E19').Value = '  +2.76%  '
$ws.Range('D20').Value = '''19.68'
$ws.Range('E20').Value = '  +2.73%  '
$ws.Range('D21').Value = '''6.973'
$ws.Range('E21').Value = '  +4.67%  '
$ws.Range('D22').Value = '''0.9994'
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').Value = '''14.18'
$ws.Range('E23').Value = '  +2.29%  '
$ws.Range('D24').Value = '24.613.89'
$ws.Range('E24').Value = '  +3.18%  '
$ws.Range('D25').Value = '''3.157'
$ws.Range('E25').Value = '  +10.37%  '
$ws.Range('E26').Value = '  +1.75%  '
$ws.Range('D27').Value = '''22.36'
$ws.Range('E27').Value = '  +2.61%  '
$ws.Range('D28').Value = '''161.19'
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '''136.45'
$ws.Range('E29').Value = '  +4.63%  '
$ws.Range('B30').Value = 'HuobiToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D30').Value = '''5.192'
$ws.Range('E30').Value = '  +1.53%  '
$ws.Range('D31').Value = '''7.529'
$ws.Range('E31').Value = '  +11.61%  '
$ws.Range('D32').Value = '1.880.42'
$ws.Range('E32').Value = '  +1.87%  '
$ws.Range('D33').Value = '''1.084'
$ws.Range('E33').Value = '  -3.42%  '
$ws.Range('D34').Value = '''0.08574'
$ws.Range('E34').Value = '  +0.66%  '
$ws.Range('D35').Value = '''7.202'
$ws.Range('E35').Value = '  +9.07%  '
$ws.Range('D36').Value = '''11.38'
$ws.Range('E36').Value = '  +8.64%  '
$ws.Range('D37').Value = '''0.2731'
$ws.Range('E37').Value = '  +3.00%  '
$ws.Range('D38').Value = '''1.936'
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('D39').Value = '''14.51'
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('D40').Value = '''0.02746'
$ws.Range('E40').Value = '  +8.49%  '
$ws.Range('D41').Value = '''0.09062'
$ws.Range('E41').Value = '  +2.93%  '
$ws.Range('D42').Value = '''1.476'
$ws.Range('E42').Value = '  +0.95%  '
$ws.Range('D43').Value = '''0.7674'
$ws.Range('E43').Value = '  +0.59%  '
$ws.Range('D44').Value = '''0.7193'
$ws.Range('E44').Value = '  +1.65%  '
$ws.Range('D45').Value = '''15.59'
$ws.Range('E45').Value = '  +2.47%  '
$ws.Range('D46').Value = '''2.535'
$ws.Range('E46').Value = '  +4.15%  '
$ws.Range('D47').Value = '''4.216'
$ws.Range('E47').Value = '  +2.45%  '
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('D49').Value = '''141.29'
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('D50').Value = '''1.322'
$ws.Range('E50').Value = '  +6.94%  '
$ws.Range('D51').Value = '''0.07997'
$ws.Range('E51').Value = '  +2.49%  '
